$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F14").Value = 'ppe'
$ws.Range("F24").Value = '135_product_information'
$ws.Range("F25").Value = 'application instructions'
$ws.Range("F26").Value = 'application instructions'
$ws.Range("F27").Value = 'env warning - water || off target movement'
$ws.Range("F29").Value = 'use restrictions'
$ws.Range("F30").Value = 'use restrictions'
$ws.Range("F31").Value = 'use restrictions'
$ws.Range("F32").Value = 'use restrictions'
$ws.Range("F35").Value = 'application instructions'
$ws.Range("F36").Value = 'application instructions'
$ws.Range("F37").Value = 'application instructions'
$ws.Range("F38").Value = 'application instructions'
$ws.Range("F39").Value = 'mixing'
$ws.Range("F40").Value = 'mixing'
$ws.Range("F73").Value = 'mixing'
$ws.Range("F75").Value = 'application instructions'
$ws.Range("F76").Value = 'mixing'
$ws.Range("F77").Value = 'mixing'
$ws.Range("F79").Value = 'application instructions'
$ws.Range("F80").Value = 'application instructions'
$ws.Range("F81").Value = 'application instructions'
$ws.Range("F82").Value = 'use restrictions'
$ws.Range("F86").Value = 'application instructions'
$ws.Range("F87").Value = 'application instructions'
$ws.Range("F88").Value = 'safety procedures'
$ws.Range("F89").Value = 'safety procedures'
$ws.Range("F90").Value = 'off target movement'
$ws.Range("F91").Value = 'off target movement'
$ws.Range("F92").Value = 'off target movement'
$ws.Range("F93").Value = 'off target movement'
$ws.Range("F94").Value = 'off target movement'
$ws.Range("F95").Value = 'off target movement'
$ws.Range("F96").Value = 'off target movement'
$ws.Range("F97").Value = 'off target movement'
$ws.Range("F98").Value = 'off target movement'
$ws.Range("F99").Value = 'application instructions'
$ws.Range("F100").Value = 'application instructions'
$ws.Range("F101").Value = '172_sensitive_areas'
$ws.Range("F102").Value = 'off target movement'
$ws.Range("F104").Value = '154_pesticide_storage'
